$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "Rounding Error:"
$ws.Range("C6").Value = 0.005
$ws.Range("C6").Style = $ws.Range("C5").Style
$ws.Range("C6").NumberFormat = $ws.Range("C5").NumberFormat

$ws.Columns.Item(2).ColumnWidth = 16

$ws.Range("B6").Select()
